$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 87: update drug name text (ESSENTIALE)
$ws.Range("A87").Value = "ESSENTIALE  "

# Append new rows 103-163
$cA = $ws.Range("A103")
$cA.Value = "ASPIRIN 81 MG ( L )"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D103").Value = "I694"

$cA = $ws.Range("A104")
$cA.Value = "ATIVAN 1 MG***SA6"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D104").Value = "G442"

$cA = $ws.Range("A105")
$cA.Value = "MYDOCALM  TAB."
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D105").Value = "G442"

$cA = $ws.Range("A106")
$cA.Value = "CAFERGOT (L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D106").Value = "G442"

$cA = $ws.Range("A107")
$cA.Value = "SYRINGE 5 CC"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D107").Value = 9929

$cA = $ws.Range("A108")
$cA.Value = "AIR - X (DISFLATYL)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D108").Value = "K30"

$cA = $ws.Range("A109")
$cA.Value = "ALUSIL"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D109").Value = "K30"

$cA = $ws.Range("A110")
$cA.Value = "BUSCOPAN 10 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D110").Value = "K30"

$cA = $ws.Range("A111")
$cA.Value = "MAGESTO"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D111").Value = "K30"

$cA = $ws.Range("A112")
$cA.Value = "BISOLVON 8 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D112").Value = "J00"

$cA = $ws.Range("A113")
$cA.Value = "DIMETAPP  TAB ( L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D113").Value = "J00"

$cA = $ws.Range("A114")
$cA.Value = "PARACETAMOL 500  MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D114").Value = "J00"

$cA = $ws.Range("A115")
$cA.Value = "MERISLON"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D115").Value = "R42"

$cA = $ws.Range("A116")
$cA.Value = "RIVOTRIL 0.5 mg"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D116").Value = "R42"

$cA = $ws.Range("A117")
$cA.Value = "ALLOPURINOL 300 MG (ZYLORIC)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Rows.Item(117).RowHeight = 45
$ws.Range("D117").Value = "M100"

$cA = $ws.Range("A118")
$cA.Value = "NAPROXEN 250 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D118").Value = "M6261"

$cA = $ws.Range("A119")
$cA.Value = "NORGESIC    (L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D119").Value = "M6261"

$cA = $ws.Range("A120")
$cA.Value = "CARDURA 2 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D120").Value = "I10"

$cA = $ws.Range("A121")
$cA.Value = "CARDURA 2 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D121").Value = "N40"

$cA = $ws.Range("A122")
$cA.Value = "OREDA"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D122").Value = "A058"

$cA = $ws.Range("A123")
$cA.Value = "IMODIUM  2 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D123").Value = "A058"

$cA = $ws.Range("A124")
$cA.Value = "BUSCOPAN 10 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D124").Value = "A058"

$cA = $ws.Range("A125")
$cA.Value = "ZOVIRAX CREAM 1 G"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D125").Value = "B029"

$cA = $ws.Range("A126")
$cA.Value = "ZOVIRAX 400 MG  "
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D126").Value = "B029"

$cA = $ws.Range("A127")
$cA.Value = "BRUFEN 400 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D127").Value = "B029"

$cA = $ws.Range("A128")
$cA.Value = "NEURONTIN 300 MG (L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D128").Value = "B029"

$cA = $ws.Range("A129")
$cA.Value = "PARACETAMOL 500  MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D129").Value = "B029"

$cA = $ws.Range("A130")
$cA.Value = "CHLORAMPHENICOL EAR DROP 10 ML"
$cA.WrapText = $true
$ws.Rows.Item(130).RowHeight = 45
$ws.Range("D130").Value = "H931"

$cA = $ws.Range("A131")
$cA.Value = "DRAMAMINE 50 MG "
$ws.Range("D131").Value = "H931"

$cA = $ws.Range("A132")
$cA.Value = "DERMOVATE CREAM (per G)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D132").Value = "L239"

$cA = $ws.Range("A133")
$cA.Value = "COUMADIN 5 MG  (HAD)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D133").Value = "I829"

$cA = $ws.Range("A134")
$cA.Value = "NAPROXEN 250 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D134").Value = "R252"

$cA = $ws.Range("A135")
$cA.Value = "NORGESIC    (L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D135").Value = "R252"

$cA = $ws.Range("A136")
$cA.Value = "NEOTICA  BALM 25 G"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D136").Value = "R252"

$cA = $ws.Range("A137")
$cA.Value = "T.T 0.5 ML เข็มที่2"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D137").Value = "Z235"

$cA = $ws.Range("A138")
$cA.Value = "DEX - OPH  EYE DROP"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D138").Value = "H000"

$cA = $ws.Range("A139")
$cA.Value = "VIGAMOX EYE DROP 5ml"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D139").Value = "H645"

$cA = $ws.Range("A140")
$cA.Value = "CLARITYNE 10 MG TAB*** SA4"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D140").Value = "R040"

$cA = $ws.Range("A141")
$cA.Value = "SPASURI  100 MG  TAB"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D141").Value = "N201"

$cA = $ws.Range("A142")
$cA.Value = "CRAVIT 500 mg (L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D142").Value = "H609"

$cA = $ws.Range("A143")
$cA.Value = "NAPROXEN 250 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D143").Value = "H609"

$cA = $ws.Range("A144")
$cA.Value = "DEX - OPH  EYE DROP"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D144").Value = "H609"

$cA = $ws.Range("A145")
$cA.Value = "DRAMAMINE 50 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D145").Value = "R42"

$cA = $ws.Range("A146")
$cA.Value = "MERISLON"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D146").Value = "R42"

$cA = $ws.Range("A147")
$cA.Value = "ALUSIL"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D147").Value = "K291"

$cA = $ws.Range("A148")
$cA.Value = "BUSCOPAN 10 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D148").Value = "K291"

$cA = $ws.Range("A149")
$cA.Value = "CIMETIDINE 400 MG TAB"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D149").Value = "K291"

$cA = $ws.Range("A150")
$cA.Value = "DIFELENE GEL"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D150").Value = "K30"

$cA = $ws.Range("A151")
$cA.Value = "AIR - X (DISFLATYL)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D151").Value = "K30"

$cA = $ws.Range("A152")
$cA.Value = "LOSEC 20 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D152").Value = "K30"

$cA = $ws.Range("A153")
$cA.Value = "MYDOCALM  TAB."
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D153").Value = "G442"

$cA = $ws.Range("A154")
$cA.Value = "ATIVAN 1 MG***SA6"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D154").Value = "G442"

$cA = $ws.Range("A155")
$cA.Value = "ENARIL 20 MG  "
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D155").Value = "I10"

$cA = $ws.Range("A156")
$cA.Value = "DEX - OPH  EYE DROP"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D156").Value = "H110"

$cA = $ws.Range("A157")
$cA.Value = "ATARAX 10 MG*** SA3"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D157").Value = "H110"

$cA = $ws.Range("A158")
$cA.Value = "DIPROSALIC 1 G CREAM (2)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D158").Value = "L239"

$cA = $ws.Range("A159")
$cA.Value = "NORFLOXACIN 400 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D159").Value = "N309"

$cA = $ws.Range("A160")
$cA.Value = "TLD (L)"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D160").Value = "B24"

$cA = $ws.Range("A161")
$cA.Value = "ATIVAN 1 MG***SA6 "
$ws.Range("D161").Value = "F411"

$cA = $ws.Range("A162")
$cA.Value = "PROCTOSEDYL  SUPPO"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D162").Value = "K649"

$cA = $ws.Range("A163")
$cA.Value = "DAFLON 500 MG"
$cA.WrapText = $true
$cA.VerticalAlignment = -4108
$ws.Range("D163").Value = "K649"

# Update selection to match final state
$ws.Range("D162:D163").Select()